$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "Save" column, matching the style of the other header cells
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column values for rows 2-13 (0/1 flags)
$saveValues = @(0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
